$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Slovakia")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Hungary"

$newSheet.Range("A9:U9").Copy() | Out-Null
$newSheet.Range("A10:U10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats ONLY
$excel.CutCopyMode = 0
Write-Output "done"
